{"js": "const pairs = [[\"538\u00d73=1614\", \"976\u00d72=1952\"], [\"722\u00d77=5054\", \"275\u00d76=1650\"], [\"828\u00d75=4140\", \"553\u00d78=4424\"], [\"790\u00d72=1580\", \"832\u00d78=6656\"], [\"193\u00d73=579\", \"899\u00d72=1798\"], [\"165\u00d77=1155\", \"737\u00d72=1474\"], [\"987\u00d78=7896\", \"769\u00d74=3076\"], [\"237\u00d79=2133\", \"548\u00d73=1644\"], [\"247\u00d76=1482\", \"339\u00d78=2712\"], [\"840\u00d77=5880\", \"438\u00d72=876\"], [\"238\u00d74=952\", \"950\u00d78=7600\"], [\"537\u00d73=1611\", \"278\u00d79=2502\"], [\"271\u00d73=813\", \"500\u00d74=2000\"], [\"919\u00d76=5514\", \"775\u00d79=6975\"], [\"651\u00d74=2604\", \"448\u00d77=3136\"], [\"507\u00d72=1014\", \"429\u00d75=2145\"], [\"681\u00d78=5448\", \"727\u00d72=1454\"], [\"857\u00d72=1714\", \"231\u00d76=1386\"], [\"877\u00d75=4385\", \"660\u00d79=5940\"], [\"801\u00d76=4806\", \"895\u00d73=2685\"], [\"541\u00d79=4869\", \"814\u00d74=3256\"], [\"659\u00d72=1318\", \"618\u00d75=3090\"], [\"275\u00d76=1650\", \"490\u00d78=3920\"], [\"558\u00d75=2790\", \"730\u00d72=1460\"], [\"435\u00d77=3045\", \"196\u00d78=1568\"]];\n\n// Process in reverse document order so that a newly-inserted value\n// that happens to match a later original value does not get\n// accidentally re-matched/re-replaced (search is done against the\n// current doc state, and some new values coincide with other rows'\n// original values).\nfor (let i = pairs.length - 1; i >= 0; i--) {\n  const [oldText, newText] = pairs[i];\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  // There may be multiple matches if an earlier (not-yet-processed) cell\n  // still holds this exact original text \"naturally\" as well as a\n  // not-yet-inserted duplicate; picking the LAST match keeps us aligned\n  // with processing order (later table cells first).\n  const item = results.items[results.items.length - 1];\n  item.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each (old, new) pair corresponds to one <w:t> cell in the practice\n# table, in document order. They are applied in REVERSE document order\n# because one new value (\"275\u00d76=1650\", row 2) happens to be identical\n# to another row's original value (row 23) later in the document; doing\n# the later row first guarantees each Find.Execute targets the intended,\n# still-original text instead of a value we just inserted.\n$pairs = @(\n    @(\"538\u00d73=1614\", \"976\u00d72=1952\"),\n    @(\"722\u00d77=5054\", \"275\u00d76=1650\"),\n    @(\"828\u00d75=4140\", \"553\u00d78=4424\"),\n    @(\"790\u00d72=1580\", \"832\u00d78=6656\"),\n    @(\"193\u00d73=579\", \"899\u00d72=1798\"),\n    @(\"165\u00d77=1155\", \"737\u00d72=1474\"),\n    @(\"987\u00d78=7896\", \"769\u00d74=3076\"),\n    @(\"237\u00d79=2133\", \"548\u00d73=1644\"),\n    @(\"247\u00d76=1482\", \"339\u00d78=2712\"),\n    @(\"840\u00d77=5880\", \"438\u00d72=876\"),\n    @(\"238\u00d74=952\", \"950\u00d78=7600\"),\n    @(\"537\u00d73=1611\", \"278\u00d79=2502\"),\n    @(\"271\u00d73=813\", \"500\u00d74=2000\"),\n    @(\"919\u00d76=5514\", \"775\u00d79=6975\"),\n    @(\"651\u00d74=2604\", \"448\u00d77=3136\"),\n    @(\"507\u00d72=1014\", \"429\u00d75=2145\"),\n    @(\"681\u00d78=5448\", \"727\u00d72=1454\"),\n    @(\"857\u00d72=1714\", \"231\u00d76=1386\"),\n    @(\"877\u00d75=4385\", \"660\u00d79=5940\"),\n    @(\"801\u00d76=4806\", \"895\u00d73=2685\"),\n    @(\"541\u00d79=4869\", \"814\u00d74=3256\"),\n    @(\"659\u00d72=1318\", \"618\u00d75=3090\"),\n    @(\"275\u00d76=1650\", \"490\u00d78=3920\"),\n    @(\"558\u00d75=2790\", \"730\u00d72=1460\"),\n    @(\"435\u00d77=3045\", \"196\u00d78=1568\")\n)\n\nfor ($i = $pairs.Count - 1; $i -ge 0; $i--) {\n    $old = $pairs[$i][0]\n    $new = $pairs[$i][1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $ok) {\n        throw \"No match found for: $old\"\n    }\n}\n"}
